$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number + the column letters/values that changed for that row,
# matching the authoritative commit diff cell-by-cell.
$updates = @(
    @{ Row=2; D='317.69'; E='3.55%'; G='17' },
    @{ Row=3; D='39.65'; E='0.96%'; G='17' },
    @{ Row=4; D='5.141'; E='0.87%'; G='17' },
    @{ Row=5; D='0.08214'; G='17' },
    @{ Row=6; D='2.031'; E='5.26%'; G='17' },
    @{ Row=7; D='8.289'; E='4.19%'; G='17' },
    @{ Row=8; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='0.9326'; E='0.05%'; G='17' },
    @{ Row=9; B='LiechtensteinCryptoassetsExchange'; C='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; D='0.1413'; E='-2.87%'; G='17' },
    @{ Row=10; B='WazirX'; C='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; D='0.1998'; E='3.68%'; G='17' },
    @{ Row=11; B='MandalaExchangeToken'; C='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D='0.09044'; E='-0.15%'; G='17' },
    @{ Row=12; B='BitrueCoin'; C='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D='0.03478'; E='-0.94%'; G='17' },
    @{ Row=13; B='BitMartToken'; C='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D='0.09801'; E='0.12%'; G='17' },
    @{ Row=14; B='BitForexToken'; C='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D='0.001404'; E='0.48%'; G='17' },
    @{ Row=15; B='TigerCash'; C='https://coinranking.com/coin/6hIn06L2+tigercash-tch'; D='0.006150'; E='4.52%'; G='17' },
    @{ Row=16; B='LEO'; C='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D='3.677'; E='-3.05%'; G='17' },
    @{ Row=17; B='GateToken'; C='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; D='4.286'; E='2.22%'; G='17' },
    @{ Row=18; D='3.195'; E='-6.43%'; G='17' },
    @{ Row=19; D='0.3472'; E='0.83%'; G='17' },
    @{ Row=20; D='0.1291'; E='-0.73%'; G='17' },
    @{ Row=21; D='4.903'; E='2.15%'; G='17' },
    @{ Row=22; D='0.2450'; E='-2.29%'; G='17' },
    @{ Row=23; D='0.04324'; E='-1.11%'; G='17' },
    @{ Row=24; D='0.001226'; E='-0.98%'; G='17' },
    @{ Row=25; E='11.54%'; G='17' },
    @{ Row=26; D='0.0001300'; E='-0.17%'; G='17' },
    @{ Row=27; D='0.0003998'; E='-10.12%'; G='17' },
    @{ Row=28; G='17' },
    @{ Row=29; G='17' },
    @{ Row=30; G='17' },
    @{ Row=31; G='17' },
    @{ Row=32; G='17' },
    @{ Row=33; G='17' },
    @{ Row=34; G='17' },
    @{ Row=35; G='17' },
    @{ Row=36; G='17' },
    @{ Row=37; G='17' },
    @{ Row=38; G='17' },
    @{ Row=39; D='0.02218'; E='8.28%'; G='17' },
    @{ Row=40; D='0.05221'; E='3.70%'; G='17' },
    @{ Row=41; D='0.007521'; E='1.18%'; G='17' },
    @{ Row=42; D='0.009785'; E='-3.25%'; G='17' },
    @{ Row=43; E='2.05%'; G='17' },
    @{ Row=44; D='0.002150'; E='0.29%'; G='17' },
    @{ Row=45; D='0.009853'; E='8.79%'; G='17' },
    @{ Row=46; D='0.00006596'; E='6.53%'; G='17' },
    @{ Row=47; D='0.00000000750'; E='-0.17%'; G='17' },
    @{ Row=48; B='BOLO'; C='https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'; D='0.002765'; E='-1.36%'; G='17' },
    @{ Row=49; B='CoinbaseStockToken'; C='https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'; D='0.001199'; E='-25.12%'; G='17' },
    @{ Row=50; D='0.00002100'; E='-0.17%'; G='17' },
    @{ Row=51; D='0.0002000'; E='-0.17%'; G='17' }
)

# Columns whose values are numeric-looking text (price/volume/hour) must be
# force-formatted as Text before assignment, otherwise Excel auto-converts
# "317.69" / "3.55%" / "17" into real numbers (losing the exact text form,
# trailing zeros, and the %-sign literal) instead of keeping them as strings.
$textCols = @("D", "E", "G")

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in $u.Keys) {
        if ($col -eq "Row") { continue }
        $addr = "$col$row"
        if ($textCols -contains $col) {
            $ws.Range($addr).NumberFormat = "@"
        }
        $ws.Range($addr).Value = $u[$col]
    }
}